$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Part 1: the "_GoBack" bookmark that used to sit between
# "...Verkeersveiligheid" and """. Verkeersveiligheid is..." is gone;
# the two runs around it are merged back into one contiguous run. A
# same-for-same Find/Replace across the old bookmark collapses it
# away and re-joins the text into a single run.
# -------------------------------------------------------------------
$d.Content.Find.Execute("sub onderwerp “Verkeersveiligheid”. Verkeersveiligheid is heel belangrijk",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "sub onderwerp “Verkeersveiligheid”. Verkeersveiligheid is heel belangrijk", 2)

# -------------------------------------------------------------------
# Part 2: spelling fix "autogarage" -> "autogarages" (adds an "s").
# The "s" ends up as its own run, immediately followed by a fresh
# "_GoBack" bookmark (Word always drops it at the most-recent edit
# point), then the remainder of the original sentence.
# -------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("autogarage", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)

$insertPos = $d.Range($rng2.End, $rng2.End)
$insertPos.InsertAfter("s")

# Force the freshly-typed "s" to live in its own run even though its
# formatting matches its neighbours (mirrors Word leaving the
# just-typed character as a distinct run until the next save/merge
# pass).
$sRange = $d.Range($rng2.End, $rng2.End + 1)
$sRange.Bold = 1
$sRange.Bold = 0

# Drop the "_GoBack" bookmark right after the newly typed "s", same
# spot Word leaves it after the most recent keystroke.
$bmPos = $d.Range($rng2.End + 1, $rng2.End + 1)
$d.Bookmarks.Add("_GoBack", $bmPos)
